$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.7044198058003014
$ws.Range("C2").Value = 0.9793093450307613
$ws.Range("D2").Value = 0.6234558220021471
$ws.Range("G2").Value = 0.4760219657335256
$ws.Range("H2").Value = 0.998

# Row 3
$ws.Range("B3").Value = 0.4033619132383097
$ws.Range("C3").Value = 0.9944254291430121
$ws.Range("D3").Value = 0.5221831119622958
$ws.Range("G3").Value = 0.4760219657335256
$ws.Range("H3").Value = 0.998

# Row 4
$ws.Range("B4").Value = 0.1541844837939589
$ws.Range("C4").Value = 0.9984108534017565
$ws.Range("D4").Value = 0.3330794258732463
$ws.Range("G4").Value = 0.4760219657335256
$ws.Range("H4").Value = 0.998

# Row 5
$ws.Range("B5").Value = 0.5261219847194402
$ws.Range("C5").Value = 0.9968660414515863
$ws.Range("D5").Value = 0.5999788466200975
$ws.Range("G5").Value = 0.4760219657335256
$ws.Range("H5").Value = 0.998
